$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the document comment text (B2): "Frame Length" -> "Frame Length / body length"
$ws.Range("B2").Value = "Frame Length / body length"

# Column B needs to widen to fit the new, longer text (bestFit-like behavior)
$ws.Range("B1:B6").EntireColumn.AutoFit() | Out-Null

# Move the active selection to B3 (matches the final saved cursor position)
$ws.Range("B3").Select()
